# Applies the AT04 Knowledge Questions (Part 2) edit:
#  1. "Concepts & Storyboards" paragraph -> switch to MyStyle, drop direct
#     run/paragraph-mark Tahoma/22 run formatting now redundant with the style.
#  2. The following hyperlink paragraph -> same pPr switch to MyStyle; the
#     hyperlink run itself keeps its existing formatting, only the trailing
#     space run and the paragraph mark's rPr lose the now-redundant
#     Tahoma/22 formatting.
#  3. The empty "Autodesk Maya" answer cell gets the pStyle MyStyle added
#     (rPr kept) plus a new run with the text "Autodesk Maya".
#  4. The empty "3ds Max" answer cell gets the same treatment with
#     "3ds Max".

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParaXml($paragraphIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $r = $p.Range
    $xml = $pkgHeader + $innerXml + $pkgFooter
    $r.InsertXML($xml)
}

# --- 1. "Concepts & Storyboards" ---
Set-ParaXml 156 '<w:p w14:paraId="2EC21FDB" w14:textId="77777777" w:rsidR="008604C6" w:rsidRDefault="00EC1899" w:rsidP="008604C6"><w:pPr><w:pStyle w:val="MyStyle"/></w:pPr><w:r><w:t>Concepts &amp; Storyboards</w:t></w:r></w:p>'

# --- 2. hyperlink paragraph ---
Set-ParaXml 157 '<w:p w14:paraId="3D29326B" w14:textId="1D9A4766" w:rsidR="00FA46B9" w:rsidRDefault="00000000" w:rsidP="008604C6"><w:pPr><w:pStyle w:val="MyStyle"/></w:pPr><w:hyperlink r:id="rId14" w:history="1"><w:r w:rsidR="00FA46B9" w:rsidRPr="00163F38"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="22"/></w:rPr><w:t>https://www.proglobalbusinesssolutions.com/3d-animation-production-process/</w:t></w:r></w:hyperlink><w:r w:rsidR="00FA46B9"><w:t xml:space="preserve"> </w:t></w:r></w:p>'

# --- 3. "Autodesk Maya" ---
Set-ParaXml 165 '<w:p w14:paraId="7A2F1143" w14:textId="0676DA24" w:rsidR="008604C6" w:rsidRDefault="008604C6" w:rsidP="008604C6"><w:pPr><w:pStyle w:val="MyStyle"/><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="22"/></w:rPr><w:t>Autodesk Maya</w:t></w:r></w:p>'

# --- 4. "3ds Max" ---
Set-ParaXml 168 '<w:p w14:paraId="3AA32B13" w14:textId="295A4DB7" w:rsidR="008604C6" w:rsidRDefault="008604C6" w:rsidP="008604C6"><w:pPr><w:pStyle w:val="MyStyle"/><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="22"/></w:rPr><w:t>3ds Max</w:t></w:r></w:p>'

Write-Host "AT04 edits applied"
